# Populate "Topic Covered" entries for the three previously-blank rows
# (dates 2024-03-14, 2024-03-15, 2024-03-16) and normalize the row
# heights of the following three rows (42-44) to match the rest of the
# sheet, per the "abstract class and methods" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C34").Value = "Virtual Methods and Abstract Class and Methods"
$ws.Range("C35").Value = "Abstract Class and Methods"
$ws.Range("C36").Value = "Holiday : Saturday"

$ws.Rows.Item(42).RowHeight = 19.5
$ws.Rows.Item(43).RowHeight = 19.5
$ws.Rows.Item(44).RowHeight = 19.5
